$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.109.23"
$ws.Range("E2").Value = "  -1.81%  "
$ws.Range("D3").Value = "1.837.04"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.66"
$ws.Range("E5").Value = "  -1.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6807"
$ws.Range("E6").Value = "  -2.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2985"
$ws.Range("E8").Value = "  -2.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07460"
$ws.Range("E9").Value = "  -3.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.13"
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07663"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").Value = "1.823.57"
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.023"
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6764"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.01"
$ws.Range("E15").Value = "  -6.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.156"
$ws.Range("E16").Value = "  -6.34%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "29.088.70"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008275"
$ws.Range("E18").Value = "  -0.57%  "
$ws.Range("D19").Value = "2.068.08"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.99"
$ws.Range("E20").Value = "  -5.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.48"
$ws.Range("E21").Value = "  -2.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.343"
$ws.Range("E23").Value = "  -3.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.72"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1439"
$ws.Range("E26").Value = "  -4.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.704"
$ws.Range("E27").Value = "  -2.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.02"
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.505"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.248"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.135"
$ws.Range("E31").Value = "  -1.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.196"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05392"
$ws.Range("E33").Value = "  +5.90%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.861"
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7496"
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.129"
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.682"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "1.301.41"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01813"
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.711"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9367"
$ws.Range("E41").Value = "  -3.77%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.080"
$ws.Range("E42").Value = "  +4.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.66"
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9985"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("B45").Value = "XinFinNetwork"
$ws.Range("C45").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.08000"
$ws.Range("E45").Value = "  +25.71%  "
$ws.Range("D46").Value = "1.970.74"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5173"
$ws.Range("E47").Value = "  -0.80%  "
$ws.Range("E48").Value = "  -3.57%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.764"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.83"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.370"
$ws.Range("E51").Value = "  -4.08%  "
